# Applies the "Enigma machine" -> "Chemistry" rewrite to the document.
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "WARNING: text not found -> $old"
    }
}

# --- Title ---
Replace-Text "Unveiling the Enigmatic Enigma Machine" "Chemistry - Unveiling the Symphony of Matter"

# --- Author name (merges 3 runs: "James C" + "." + " Ellis") ---
Replace-Text "James C. Ellis" "Harper Anderson"

# --- Email (keeps the separating "." run untouched, merges the rest) ---
Replace-Text "ellis" "harperanderson03@highschool"
Replace-Text "james@cryptography.research" "edu"

# --- Body paragraph (first block, separated by <w:br/> line breaks) ---
Replace-Text "From the labyrinthine depths of human intelligence emerged a mechanical marvel that would forever alter the course of secrecy and warfare: the Enigma machine" "Chemistry, an enthralling realm where matter transforms and elements dance in a mesmerizing ballet, beckons us to unlock the secrets of our physical world"

Replace-Text " This intricate electromechanical device, born in the 1920s, became the linchpin of German communications during the Second World War, shrouding messages in a bewildering tapestry of encryptions that defied comprehension" " It is a field of perpetual wonder, revealing the complexities of life and the universe we inhabit"

Replace-Text " The Enigma machine, with its enigmatic scrambling of letters, challenged the world's most brilliant codebreakers, propelling them into a relentless race against time to unravel its impenetrable secrets" " As we embark on this captivating journey, we will explore the elements, the building blocks of all things, delving into their properties and behaviors. We will unravel the intricacies of chemical bonds, understanding how they hold atoms together in a symphony of molecular architectures. Furthermore, we will uncover the dynamics of chemical reactions, witnessing the energy transformations that shape our world"

# --- Second block ---
Replace-Text "In a clandestine world where nations fought a silent battle of wits, the Enigma machine stood as an impenetrable fortress, guarding military strategies and diplomatic maneuvers from prying eyes" "In the realm of chemistry, we will investigate the fascinating world of compounds, where elements unite in harmonious combinations, displaying a vast array of unique characteristics"

Replace-Text " Each keystroke on its intricate keyboard triggered a mesmerizing dance of rotating rotors, scrambling the alphabet in a seemingly random, indecipherable manner" " Through the lens of chemistry, we will decipher the mysteries of acids and bases, exploring their reactivity and impact on our daily lives"

Replace-Text " The resulting messages, resembling cryptic puzzles, tantalized codebreakers, who tirelessly sought to penetrate this formidable barrier" " We will delve into the intricacies of organic chemistry, discovering the remarkable versatility of carbon-containing molecules and their significance in living organisms. Moreover, we will unveil the marvels of biochemistry, unraveling the intricate mechanisms that govern cellular processes and sustain life"

# --- Third block ---
Replace-Text "As the world plunged into the abyss of global conflict, the Enigma machine became an indispensable tool for military communication" "As we continue our exploration, we will unravel the profound impact chemistry has on our lives and the world around us"

Replace-Text " German commanders relied heavily on its impenetrable encryption to coordinate troop movements, convey tactical plans, and share intelligence reports" " From the marvels of modern medicine to the wonders of materials science, we will witness the practical applications of chemistry in diverse fields"

Replace-Text " The stakes were extraordinarily high, as the fate of nations hung precariously upon the ability to safeguard sensitive information fromDi Fang  intelligence efforts" " Through the study of chemistry, we will gain a deeper appreciation for the intricate workings of the natural world, empowering us to make informed decisions and solve real-world problems. We will emerge with a profound understanding of the chemical world, enabling us to contribute to the betterment of society and tackle global challenges with innovative solutions"

# --- Summary paragraph ---
Replace-Text "The Enigma machine, an ingenious electromechanical device, played a pivotal role in German communication during World War II" "Chemistry, a captivating science, unravels the symphony of matter and propels us into the depths of the physical world"

Replace-Text " Its intricate system of rotating rotors and complex encryptions baffled codebreakers, making it a formidable challenge to decipher" " We delve into the realm of elements, compounds, and chemical reactions, deciphering the intricate dance of molecules and understanding the properties and behaviors that govern their interactions"

Replace-Text " The Enigma machine became a symbol of German military prowess, safeguarding sensitive information and enabling effective coordination among their forces" " Chemistry offers a profound perspective on life and the universe, empowering us to appreciate the wonders of our surroundings and apply our knowledge to solve real-world problems"

Replace-Text " However, the determined efforts of Allied codebreakers, notably Alan Turing and his team at Bletchley Park, ultimately cracked the Enigma code, turning the tide of the war in favor of the Allies" " It is a testament to the interconnectedness of all things and the boundless possibilities that lie within the realm of scientific discovery"

Replace-Text " The Enigma machine remains a testament to human ingenuity and the relentless pursuit of deciphering even the most enigmatic secrets" " Through chemistry, we unlock the secrets of the universe and forge a path toward a brighter future for humanity"

# --- Add a new empty paragraph at the very end of the document ---
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
